# Generate Report for Handoff
#
# For the files that are "Ready for handoff" but have not yet received a
# handback file (the "Latest Handback File" column is still empty), flag
# them with handoff-type priority "ht" in the localization status report,
# and refresh the handoff-generation timestamps for those same rows on the
# Overview sheet and on each locale sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Rows (on each locale sheet) whose handback hasn't arrived yet.
$rows = @(7, 9, 11, 12, 13, 14)

foreach ($r in $rows) {
    # Mark the Priority column with the handoff type.
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"

    # Refresh the "Latest Handoff Datetime" on each locale sheet.
    $wsZhCn.Range("H$r").Value = "2016-08-31 14:26:01"
    $wsDeDe.Range("H$r").Value = "2016-08-31 14:26:18"

    # Refresh the "Latest HO Xliff Generate Date" on the Overview sheet.
    $wsOverview.Range("G$r").Value = "2016-08-31 14:26:18"
}
